# "mise a jour bouton engagement htp"
#
# - Feuil12: move the active selection to A18, widen column B (33 chars)
#   and give column I an explicit width.
# - Feuil11 is renamed to "Feuil3", its first column is widened, a third
#   data row (identical in shape to rows 1-2, pointing at htpsalesstocktot)
#   is appended together with its hyperlink, and it becomes the active
#   (selected) sheet of the workbook instead of Feuil1.

$wb = $excel.ActiveWorkbook

# ---- Feuil12 (sheet4): selection + column widths -------------------------
$ws4 = $wb.Worksheets.Item("Feuil12")
$ws4.Range("A18").Select()
$ws4.Columns.Item(2).ColumnWidth = 32.15   # -> stored width 33
$ws4.Columns.Item(9).ColumnWidth = 28.0    # -> stored width ~28.8555 (closest reachable)

# ---- Feuil11 -> Feuil3 (sheet3) -------------------------------------------
$ws3 = $wb.Worksheets.Item("Feuil11")
$ws3.Name = "Feuil3"
$ws3.Columns.Item(1).ColumnWidth = 98.3    # -> stored width ~99.1406 (closest reachable)

# New row 3, same layout as rows 1-2 but referencing "htpsalesstocktot"
$ws3.Range("A3").Value = "//10.128.1.2/bpo_almerys/00-TOUS/06-DOSSIER POLE/01-HTP/05- REPORTING/03-HTP/DOC_HTP"
$ws3.Range("B3").Value = "/SALESFORCE/"
$ws3.Range("C3").Value = "DATE"
$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = "Ligne"
$ws3.Range("F3").Value = "htpsalesstocktot"
$ws3.Range("G3").Value = "Ligne"
$ws3.Range("H3").Value = 0
$ws3.Range("I3").Value = "SALESFORCE/"
$ws3.Range("J3").Value = "TRAME SALESFORCE"

# Hyperlink for the new A3 cell, matching the style used by A1/A2
$ws3.Hyperlinks.Add($ws3.Range("A3"), "file:///\\10.128.1.2\bpo_almerys\", $null, $null, "\\10.128.1.2\bpo_almerys\")
$ws3.Range("A3").Value = "//10.128.1.2/bpo_almerys/00-TOUS/06-DOSSIER POLE/01-HTP/05- REPORTING/03-HTP/DOC_HTP"
$ws3.Range("A1").Copy()
$ws3.Range("A3").PasteSpecial(-4122)   # xlPasteFormats: reuse A1's "Lien hypertexte" style

# Make Feuil3 the active sheet/tab with its selection on A9
$ws3.Range("A9").Select()
